# Update cryptos list: apply price/volume/name/link changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '58.484.14'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.624.79'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'535.19"
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').Value = "'142.92"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = "'6.97"
$ws.Range('E9').Value = '  +7.28%  '
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').Value = '3.091.33'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').Value = '58.428.11'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').Value = "'20.75"
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '2.649.61'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('D19').Value = "'334.65"
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = "'66.23"
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').Value = '0.0₃0736'
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').Value = "'18.75"
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = "'150.26"
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('D35').Value = "'37.17"
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = "'0.849"
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'1.10"
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = "'1.41"
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = "'0.809"
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D41').Value = "'280.15"
$ws.Range('E41').Value = '  +2.74%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = "'10.67"
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = "'0.0530"
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'19.00"
$ws.Range('E46').Value = '  +3.01%  '
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').Value = '1.945.52'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').Value = "'4.44"
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -3.90%  '
